# Updated symbol list on Wed Jan 18 13:11:53 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) / Hora columns for each coin row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal text value (matching the source
    # inline-string cells) instead of Excel auto-coercing numeric- or
    # percent-looking strings into numbers, then drop the temporary
    # "@" text format so the cell keeps its original (default) style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '301.60'
Set-TextValue $ws.Range("E2") '-0.08%'
Set-TextValue $ws.Range("G2") '13'
Set-TextValue $ws.Range("D3") '32.31'
Set-TextValue $ws.Range("E3") '1.49%'
Set-TextValue $ws.Range("G3") '13'
Set-TextValue $ws.Range("D4") '5.054'
Set-TextValue $ws.Range("E4") '-1.39%'
Set-TextValue $ws.Range("G4") '13'
Set-TextValue $ws.Range("D5") '0.07680'
Set-TextValue $ws.Range("E5") '-2.02%'
Set-TextValue $ws.Range("G5") '13'
Set-TextValue $ws.Range("D6") '2.065'
Set-TextValue $ws.Range("E6") '-8.97%'
Set-TextValue $ws.Range("G6") '13'
Set-TextValue $ws.Range("D7") '7.851'
Set-TextValue $ws.Range("E7") '0.62%'
Set-TextValue $ws.Range("G7") '13'
Set-TextValue $ws.Range("D8") '3.773'
Set-TextValue $ws.Range("E8") '-1.08%'
Set-TextValue $ws.Range("G8") '13'
Set-TextValue $ws.Range("D9") '0.9193'
Set-TextValue $ws.Range("E9") '-0.90%'
Set-TextValue $ws.Range("G9") '13'
Set-TextValue $ws.Range("D10") '0.1763'
Set-TextValue $ws.Range("E10") '-0.50%'
Set-TextValue $ws.Range("G10") '13'
Set-TextValue $ws.Range("D11") '0.07865'
Set-TextValue $ws.Range("E11") '3.24%'
Set-TextValue $ws.Range("G11") '13'
Set-TextValue $ws.Range("D12") '0.08441'
Set-TextValue $ws.Range("E12") '-4.94%'
Set-TextValue $ws.Range("G12") '13'
Set-TextValue $ws.Range("D13") '0.03062'
Set-TextValue $ws.Range("E13") '-1.14%'
Set-TextValue $ws.Range("G13") '13'
Set-TextValue $ws.Range("D14") '0.09975'
Set-TextValue $ws.Range("E14") '-0.47%'
Set-TextValue $ws.Range("G14") '13'
Set-TextValue $ws.Range("D15") '0.001515'
Set-TextValue $ws.Range("E15") '0.28%'
Set-TextValue $ws.Range("G15") '13'
Set-TextValue $ws.Range("D16") '0.005898'
Set-TextValue $ws.Range("E16") '-1.75%'
Set-TextValue $ws.Range("G16") '13'
Set-TextValue $ws.Range("D17") '0.007498'
Set-TextValue $ws.Range("E17") '2,116.77%'
Set-TextValue $ws.Range("G17") '13'
Set-TextValue $ws.Range("D18") '3.470'
Set-TextValue $ws.Range("E18") '-0.13%'
Set-TextValue $ws.Range("G18") '13'
Set-TextValue $ws.Range("D19") '2.151'
Set-TextValue $ws.Range("E19") '-4.44%'
Set-TextValue $ws.Range("G19") '13'
Set-TextValue $ws.Range("D20") '0.3339'
Set-TextValue $ws.Range("E20") '1.42%'
Set-TextValue $ws.Range("G20") '13'
Set-TextValue $ws.Range("D21") '0.1324'
Set-TextValue $ws.Range("E21") '-1.04%'
Set-TextValue $ws.Range("G21") '13'
Set-TextValue $ws.Range("D22") '4.266'
Set-TextValue $ws.Range("E22") '-1.42%'
Set-TextValue $ws.Range("G22") '13'
Set-TextValue $ws.Range("D23") '0.1974'
Set-TextValue $ws.Range("E23") '10.10%'
Set-TextValue $ws.Range("G23") '13'
Set-TextValue $ws.Range("D24") '0.04530'
Set-TextValue $ws.Range("E24") '-1.76%'
Set-TextValue $ws.Range("G24") '13'
Set-TextValue $ws.Range("D25") '0.001230'
Set-TextValue $ws.Range("E25") '-1.87%'
Set-TextValue $ws.Range("G25") '13'
Set-TextValue $ws.Range("D26") '0.004126'
Set-TextValue $ws.Range("E26") '-7.84%'
Set-TextValue $ws.Range("G26") '13'
Set-TextValue $ws.Range("D27") '0.0001250'
Set-TextValue $ws.Range("E27") '-0.01%'
Set-TextValue $ws.Range("G27") '13'
Set-TextValue $ws.Range("D28") '--'
Set-TextValue $ws.Range("E28") '--%'
Set-TextValue $ws.Range("G28") '13'
Set-TextValue $ws.Range("D29") '--'
Set-TextValue $ws.Range("E29") '--%'
Set-TextValue $ws.Range("G29") '13'
Set-TextValue $ws.Range("D30") '--'
Set-TextValue $ws.Range("E30") '--%'
Set-TextValue $ws.Range("G30") '13'
Set-TextValue $ws.Range("D31") '--'
Set-TextValue $ws.Range("E31") '--%'
Set-TextValue $ws.Range("G31") '13'
Set-TextValue $ws.Range("D32") '--'
Set-TextValue $ws.Range("E32") '--%'
Set-TextValue $ws.Range("G32") '13'
Set-TextValue $ws.Range("D33") '--'
Set-TextValue $ws.Range("E33") '--%'
Set-TextValue $ws.Range("G33") '13'
Set-TextValue $ws.Range("D34") '--'
Set-TextValue $ws.Range("E34") '--%'
Set-TextValue $ws.Range("G34") '13'
Set-TextValue $ws.Range("D35") '--'
Set-TextValue $ws.Range("E35") '--%'
Set-TextValue $ws.Range("G35") '13'
Set-TextValue $ws.Range("D36") '--'
Set-TextValue $ws.Range("E36") '--%'
Set-TextValue $ws.Range("G36") '13'
Set-TextValue $ws.Range("D37") '--'
Set-TextValue $ws.Range("E37") '--%'
Set-TextValue $ws.Range("G37") '13'
Set-TextValue $ws.Range("D38") '--'
Set-TextValue $ws.Range("E38") '--%'
Set-TextValue $ws.Range("G38") '13'
Set-TextValue $ws.Range("D39") '0.01714'
Set-TextValue $ws.Range("E39") '-3.32%'
Set-TextValue $ws.Range("G39") '13'
Set-TextValue $ws.Range("D40") '0.04674'
Set-TextValue $ws.Range("E40") '-1.99%'
Set-TextValue $ws.Range("G40") '13'
Set-TextValue $ws.Range("D41") '0.007510'
Set-TextValue $ws.Range("E41") '2.34%'
Set-TextValue $ws.Range("G41") '13'
Set-TextValue $ws.Range("D42") '0.1350'
Set-TextValue $ws.Range("E42") '-1.02%'
Set-TextValue $ws.Range("G42") '13'
Set-TextValue $ws.Range("D43") '0.002329'
Set-TextValue $ws.Range("E43") '6.39%'
Set-TextValue $ws.Range("G43") '13'
Set-TextValue $ws.Range("D44") '0.01043'
Set-TextValue $ws.Range("E44") '7.39%'
Set-TextValue $ws.Range("G44") '13'
Set-TextValue $ws.Range("D45") '0.00006229'
Set-TextValue $ws.Range("E45") '-0.65%'
Set-TextValue $ws.Range("G45") '13'
Set-TextValue $ws.Range("D46") '0.00000000749'
Set-TextValue $ws.Range("E46") '-0.11%'
Set-TextValue $ws.Range("G46") '13'
Set-TextValue $ws.Range("D47") '0.002998'
Set-TextValue $ws.Range("E47") '-62.47%'
Set-TextValue $ws.Range("G47") '13'
Set-TextValue $ws.Range("D48") '0.8206'
Set-TextValue $ws.Range("E48") '16.04%'
Set-TextValue $ws.Range("G48") '13'
Set-TextValue $ws.Range("D49") '0.00002099'
Set-TextValue $ws.Range("E49") '-0.11%'
Set-TextValue $ws.Range("G49") '13'
Set-TextValue $ws.Range("D50") '0.0001999'
Set-TextValue $ws.Range("E50") '-0.11%'
Set-TextValue $ws.Range("G50") '13'
Set-TextValue $ws.Range("D51") '--'
Set-TextValue $ws.Range("E51") '--%'
Set-TextValue $ws.Range("G51") '13'
